# Update the panel_query_time-stamped "time_taken" values on the "data" sheet,
# then add a new "metadata" sheet (placed after "data") describing the panel query.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Refresh the per-row query timestamps in column F of the "data" sheet.
$ws.Range("F2").Value = "2021-10-05 14:35:17.276100"
$ws.Range("F3").Value = "2021-10-05 14:35:17.276108"
$ws.Range("F4").Value = "2021-10-05 14:35:17.276111"
$ws.Range("F5").Value = "2021-10-05 14:35:17.276114"
$ws.Range("F6").Value = "2021-10-05 14:35:17.276117"
$ws.Range("F7").Value = "2021-10-05 14:35:17.276120"
$ws.Range("F8").Value = "2021-10-05 14:35:17.276122"
$ws.Range("F9").Value = "2021-10-05 14:35:17.276125"
$ws.Range("F10").Value = "2021-10-05 14:35:17.276128"
$ws.Range("F11").Value = "2021-10-05 14:35:17.276130"

# 2) Add the new "metadata" worksheet right after "data".
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Copy the header formatting (bold, centered, bordered) from "data" so the new
# sheet's styled cells reuse the same cell style (s="1").
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Porphyria"
$meta.Range("C2").Value = 3077

# "data_version" must stay text ("0.12"), not be coerced to the number 0.12.
# Force text entry via NumberFormat, then restore the plain (unstyled) look by
# re-pasting (formats only) the default formatting from an untouched, unstyled
# cell -- this changes the visual style only, the stored text value is kept.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.12"
$ws.Range("C2").Copy()
$meta.Range("D2").PasteSpecial(-4122)

$meta.Range("E2").Value = "2020-07-22T02:36:09.647302Z"
$meta.Range("F2").Value = "2021-10-05 14:35:17.272367"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3077/?format=json"

$ws.Select()
$ws.Range("A1").Select()
